# Apply scheduled-runner price/profit updates to Sheets/Hades_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2027.7037
$ws.Range("I62").Value = 1838.55
$ws.Range("J62").Value = 2568.1428
$ws.Range("K62").Value = 1838.55
$ws.Range("L62").Value = 2568.1428
$ws.Range("M62").Value = -1214.55
$ws.Range("N62").Value = -3816.1428
$ws.Range("H65").Value = 2027.7037
$ws.Range("I65").Value = 1838.55
$ws.Range("J65").Value = 2568.1428
$ws.Range("K65").Value = 9192.75
$ws.Range("L65").Value = 12840.714
$ws.Range("M65").Value = -6072.75
$ws.Range("N65").Value = -19080.714
$ws.Range("H125").Value = 608.64703
$ws.Range("I125").Value = 453.9091
$ws.Range("J125").Value = 892.3333
$ws.Range("K125").Value = 4085.1819
$ws.Range("L125").Value = 8030.9997
$ws.Range("M125").Value = -1625.1819
$ws.Range("N125").Value = -12950.9997
$ws.Range("H135").Value = 50845.273
$ws.Range("I135").Value = 40969.76
$ws.Range("J135").Value = 67304.47
$ws.Range("K135").Value = 368727.84
$ws.Range("L135").Value = 605740.23
$ws.Range("M135").Value = -366192.84
$ws.Range("N135").Value = -610810.23
$ws.Range("H136").Value = 49105
$ws.Range("J136").Value = 49105
$ws.Range("L136").Value = 49105
$ws.Range("N136").Value = -59305
$ws.Range("H137").Value = 2858953.8
$ws.Range("I137").Value = 5264779.5
$ws.Range("K137").Value = 15794338.5
$ws.Range("M137").Value = -15791788.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8085.3335
$ws.Range("I31").Value = 8085.3335
$ws.Range("K31").Value = 8085.3335
$ws.Range("M31").Value = -7791.3335
$ws.Range("H32").Value = 5396401
$ws.Range("I32").Value = 5868591.5
$ws.Range("J32").Value = 13428.2
$ws.Range("K32").Value = 5868591.5
$ws.Range("L32").Value = 13428.2
$ws.Range("M32").Value = -5868304.5
$ws.Range("N32").Value = -14002.2
$ws.Range("H45").Value = 2708.5173
$ws.Range("J45").Value = 1763.75
$ws.Range("L45").Value = 1763.75
$ws.Range("N45").Value = -2517.75
$ws.Range("H61").Value = 35786636
$ws.Range("I61").Value = 55612076
$ws.Range("J61").Value = 100842.8
$ws.Range("K61").Value = 55612076
$ws.Range("L61").Value = 100842.8
$ws.Range("M61").Value = -55611864
$ws.Range("N61").Value = -101266.8
$ws.Range("H74").Value = 7412710.5
$ws.Range("I74").Value = 14765453
$ws.Range("J74").Value = 59967.65
$ws.Range("K74").Value = 14765453
$ws.Range("L74").Value = 59967.65
$ws.Range("M74").Value = -14764579
$ws.Range("N74").Value = -61715.65
$ws.Range("H77").Value = 7412710.5
$ws.Range("I77").Value = 14765453
$ws.Range("J77").Value = 59967.65
$ws.Range("K77").Value = 73827265
$ws.Range("L77").Value = 299838.25
$ws.Range("M77").Value = -73822897
$ws.Range("N77").Value = -308574.25
$ws.Range("H122").Value = 4832977.5
$ws.Range("I122").Value = 2098.5789
$ws.Range("J122").Value = 27779652
$ws.Range("K122").Value = 6295.736699999999
$ws.Range("L122").Value = 83338956
$ws.Range("M122").Value = -3845.736699999999
$ws.Range("N122").Value = -83343856
$ws.Range("H132").Value = 73560.75
$ws.Range("I132").Value = 60625.59
$ws.Range("J132").Value = 93551.45
$ws.Range("K132").Value = 181876.77
$ws.Range("L132").Value = 280654.35
$ws.Range("M132").Value = -179346.77
$ws.Range("N132").Value = -285714.35
$ws.Range("H135").Value = 46349.875
$ws.Range("J135").Value = 46349.875
$ws.Range("L135").Value = 46349.875
$ws.Range("N135").Value = -56489.875
$ws.Range("H136").Value = 35786636
$ws.Range("I136").Value = 55612076
$ws.Range("J136").Value = 100842.8
$ws.Range("K136").Value = 166836228
$ws.Range("L136").Value = 302528.4
$ws.Range("M136").Value = -166833678
$ws.Range("N136").Value = -307628.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2040.421
$ws.Range("I107").Value = 1785.4445
$ws.Range("J107").Value = 2269.9
$ws.Range("K107").Value = 1785.4445
$ws.Range("L107").Value = 2269.9
$ws.Range("M107").Value = 134.5554999999999
$ws.Range("N107").Value = -6109.9
$ws.Range("H140").Value = 55213.684
$ws.Range("J140").Value = 55213.684
$ws.Range("L140").Value = 55213.684
$ws.Range("N140").Value = -65573.68400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 37038668
$ws.Range("I16").Value = 2041.3334
$ws.Range("J16").Value = 111111920
$ws.Range("K16").Value = 2041.3334
$ws.Range("L16").Value = 111111920
$ws.Range("M16").Value = -1754.3334
$ws.Range("N16").Value = -111112494
$ws.Range("H31").Value = 1867.1177
$ws.Range("I31").Value = 1081.5264
$ws.Range("J31").Value = 2333.5625
$ws.Range("K31").Value = 1081.5264
$ws.Range("L31").Value = 2333.5625
$ws.Range("M31").Value = -786.5264
$ws.Range("N31").Value = -2923.5625
$ws.Range("H34").Value = 1867.1177
$ws.Range("I34").Value = 1081.5264
$ws.Range("J34").Value = 2333.5625
$ws.Range("K34").Value = 1081.5264
$ws.Range("L34").Value = 2333.5625
$ws.Range("M34").Value = -879.5264
$ws.Range("N34").Value = -2737.5625
$ws.Range("H52").Value = 44060
$ws.Range("J52").Value = 44060
$ws.Range("L52").Value = 44060
$ws.Range("N52").Value = -44648
$ws.Range("H113").Value = 37038668
$ws.Range("I113").Value = 2041.3334
$ws.Range("J113").Value = 111111920
$ws.Range("K113").Value = 2041.3334
$ws.Range("L113").Value = 111111920
$ws.Range("M113").Value = 128.6666
$ws.Range("N113").Value = -111116260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1264.5
$ws.Range("I68").Value = 662.1
$ws.Range("J68").Value = 1551.3572
$ws.Range("K68").Value = 1986.3
$ws.Range("L68").Value = 4654.071599999999
$ws.Range("M68").Value = -1175.3
$ws.Range("N68").Value = -6276.071599999999
$ws.Range("H71").Value = 1264.5
$ws.Range("I71").Value = 662.1
$ws.Range("J71").Value = 1551.3572
$ws.Range("K71").Value = 5958.900000000001
$ws.Range("L71").Value = 13962.2148
$ws.Range("M71").Value = -1902.900000000001
$ws.Range("N71").Value = -22074.2148
$ws.Range("H80").Value = 3437.1667
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 3951.4614
$ws.Range("K80").Value = 6300
$ws.Range("L80").Value = 11854.3842
$ws.Range("M80").Value = -5364
$ws.Range("N80").Value = -13726.3842
$ws.Range("H83").Value = 3437.1667
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 3951.4614
$ws.Range("K83").Value = 18900
$ws.Range("L83").Value = 35563.1526
$ws.Range("M83").Value = -14220
$ws.Range("N83").Value = -44923.1526
$ws.Range("H107").Value = 1115.6198
$ws.Range("I107").Value = 481.88095
$ws.Range("J107").Value = 2033.4482
$ws.Range("K107").Value = 1445.64285
$ws.Range("L107").Value = 6100.3446
$ws.Range("M107").Value = 474.35715
$ws.Range("N107").Value = -9940.3446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4200
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4200
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 12600
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -17540
$ws.Range("H132").Value = 47863.75
$ws.Range("I132").Value = 34274.484
$ws.Range("J132").Value = 80268.92
$ws.Range("K132").Value = 102823.452
$ws.Range("L132").Value = 240806.76
$ws.Range("M132").Value = -100293.452
$ws.Range("N132").Value = -245866.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 26006.857
$ws.Range("I132").Value = 1465.0385
$ws.Range("J132").Value = 65887.31
$ws.Range("K132").Value = 4395.1155
$ws.Range("L132").Value = 197661.93
$ws.Range("M132").Value = -1865.1155
$ws.Range("N132").Value = -202721.93
$ws.Range("H137").Value = 25000
$ws.Range("J137").Value = 25000
$ws.Range("L137").Value = 25000
$ws.Range("N137").Value = -35200
$ws.Range("H139").Value = 51277.5
$ws.Range("J139").Value = 51277.5
$ws.Range("L139").Value = 51277.5
$ws.Range("N139").Value = -61557.5
$ws.Range("H141").Value = 51715
$ws.Range("J141").Value = 51715
$ws.Range("L141").Value = 51715
$ws.Range("N141").Value = -62075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 31000.586
$ws.Range("J76").Value = 31000.586
$ws.Range("L76").Value = 31000.586
$ws.Range("N76").Value = -31630.586
$ws.Range("H79").Value = 31000.586
$ws.Range("J79").Value = 31000.586
$ws.Range("L79").Value = 31000.586
$ws.Range("N79").Value = -33184.586
